# Insert a new weekly price record as row 36 (Arándano (blue), Vega Central
# Mapocho de Santiago), pushing the existing rows 36-150 down to 37-151.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(36).Insert()

$ws.Range("A36").Value = 9
$ws.Range("B36").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 44560
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100101
$ws.Range("H36").Value = "Berries"
$ws.Range("I36").Value = 100101001
$ws.Range("J36").Value = "Arándano (blue)"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 580
$ws.Range("N36").Value = 3000
$ws.Range("O36").Value = 3200
$ws.Range("P36").Value = 3097
$ws.Range("Q36").Value = "$/bandeja 2 kilos"
$ws.Range("R36").Value = "Región de O'Higgins"
$ws.Range("S36").Value = 1548
$ws.Range("T36").Value = 2
